$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 413, shifting existing rows 413:426 down to 414:427
$ws.Rows.Item(413).Insert()

# Populate the newly inserted row 413 with the new record
$ws.Cells.Item(413, 1).Value = 4
$ws.Cells.Item(413, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(413, 3).Value = "Los Lagos"
$ws.Cells.Item(413, 4).Value = 45041
$ws.Cells.Item(413, 4).NumberFormat = $ws.Cells.Item(414, 4).NumberFormat
$ws.Cells.Item(413, 5).Value = 10
$ws.Cells.Item(413, 6).Value = "Fruta"
$ws.Cells.Item(413, 7).Value = 100102
$ws.Cells.Item(413, 8).Value = "Cítricos"
$ws.Cells.Item(413, 9).Value = 100102004
$ws.Cells.Item(413, 10).Value = "Mandarina"
$ws.Cells.Item(413, 11).Value = "Murcott"
$ws.Cells.Item(413, 12).Value = "Segunda"
$ws.Cells.Item(413, 13).Value = 300
$ws.Cells.Item(413, 14).Value = 17000
$ws.Cells.Item(413, 15).Value = 17000
$ws.Cells.Item(413, 16).Value = 17000
$ws.Cells.Item(413, 17).Value = "`$/caja 12 kilos granel"
$ws.Cells.Item(413, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(413, 19).Value = 17000
$ws.Cells.Item(413, 20).Value = 1
